# Add scenarios for PDE file
# Adds two new test-case rows (4 and 5) to the MFP sheet describing a
# new scenario ("Test NDC-11 and Effective Date Dynamicly"), widens the
# sql_query / expected_value columns to fit the new content, wraps +
# top-aligns the JSON sample cell, and flips the sheet to portrait
# print orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new content -----------------------------------------------------
$query       = " select ndc_eleven, effective_date from mfp;"
$getFromFile = " getFromFile(""ndc_eleven"", ""effective_date"")"
$title       = "Test NDC-11 and Effective Date Dynamicly"
$json        = "{`n`t{ `n`t`t""ndc_eleven"": ""00003-0893-21"",`n`t`t""effective-date"": ""1/1/2026""`n`t},`n`t{`n`t`t""ndc_eleven"": ""00003-0893-31"",`n`t`t""effective_date"": ""1/1/2026""`n`t}`n}"

# Write cell-by-cell in the order that reproduces the target shared
# string table ordering (query, getFromFile, title, json).
$ws.Range("A5").Value = 4
$ws.Range("C5").Value = $query
$ws.Range("C6").Value = $query
$ws.Range("D6").Value = $getFromFile
$ws.Range("B5").Value = $title
$ws.Range("B6").Value = $title
$ws.Range("D5").Value = $json
$ws.Range("A6").Value = 5

# --- formatting for the new JSON sample cell --------------------------
$ws.Range("D5").WrapText = $true
$ws.Range("D5").VerticalAlignment = -4160  # xlTop

# Row 5 grows tall enough to show the full JSON sample.
$ws.Rows.Item(5).RowHeight = 256.2

# --- column widths: widen sql_query (C) and expected_value (D) --------
$ws.Columns.Item(3).ColumnWidth = 37.25
$ws.Columns.Item(4).ColumnWidth = 74.25

# --- selection + print setup ------------------------------------------
$ws.Range("D5").Select()
$ws.PageSetup.Orientation = 1  # xlPortrait
